# Fruta / hortaliza, semanal
# Insert a new weekly record as row 7, pushing the existing rows 7-17 down to 8-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 7 (existing rows shift down one).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly observation.
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 45014
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100104
$ws.Cells.Item(7, 8).Value = "Frutos de pepita"
$ws.Cells.Item(7, 9).Value = 100104003
$ws.Cells.Item(7, 10).Value = "Membrillo"
$ws.Cells.Item(7, 11).Value = "Champion"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 9000
$ws.Cells.Item(7, 15).Value = 10000
$ws.Cells.Item(7, 16).Value = 9500
$ws.Cells.Item(7, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(7, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(7, 19).Value = 528
$ws.Cells.Item(7, 20).Value = 18
